$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "501+ ЛЕНТА МАСКИР. ОДНОСТ., ОСНОВА БУМ., АДГЕЗИВ КАУЧУК., БЕЖ.; 0048 ММХ 0055,0 М"
$ws.Range("C31").Value = 6.7
$ws.Range("D31").Value = "нескладской / срок поставки 62 дн."
